# Rebuild the "Profiling" sheet with a cProfile-style brute-force N-body
# benchmark table (pasted from the project's profiling output), matching
# the target commit "remove h call brutes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Profiling".
$ws.Name = "Profiling"

# --- Header row (row 3) ----------------------------------------------------
$headers = @("% Time", "Cumulative Seconds", "Self Seconds", "Calls", "Self s/call", "Total s/call", "Function Name")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $headers[$i]
}
$headerRange = $ws.Range("A3:G3")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 10.5
$headerRange.Font.Color = 5325111   # RGB(55,65,81) = #374151, BGR-packed for COM

# --- Data rows (rows 4-12) --------------------------------------------------
$data = @(
    @(89.06, 3.78, 3.78, 40000,     0, 0,    "compute_brute_force"),
    @(8.73,  4.16, 0.37, 376713688, 0, 0,    "max"),
    @(1.18,  4.21, 0.05, "-",       "-", "-", "min"),
    @(0,     4.21, 0,    4,         0, 0,    "update_positions"),
    @(0,     4.21, 0,    2,         0, 0,    "second"),
    @(0,     4.21, 0,    1,         0, 0,    "get_nbr_particles"),
    @(0,     4.21, 0,    1,         0, 4.16, "nbodybruteforce"),
    @(0,     4.21, 0,    1,         0, 0,    "print_parameters"),
    @(0,     4.21, 0,    1,         0, 0,    "read_test_case")
)

$bodyRange = $ws.Range("A4:G12")
$row = 4
foreach ($r in $data) {
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $r[$col - 1]
    }
    $row++
}
$bodyRange.Font.Name = "Arial"
$bodyRange.Font.Size = 10.5
$bodyRange.Font.Color = 5325111

# "Calls" column (D) on the first two data rows uses a thousands-separated
# integer number format.
$ws.Range("D4:D5").NumberFormat = "#,##0"

# --- Title row (written last so it lands at the end of the shared-string
# table, matching how the sheet was actually authored) ---------------------
$ws.Range("A1").Value = "brute force"

# --- View tweaks -------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 181
$ws.Range("A16").Select()
